$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("I1").Value = "SiteID"
$ws.Range("J1").Value = "Location ID"
$ws.Range("K1").Value = "Location Number"

# Row 2 data
$ws.Range("H2").Value = "a811K0000004fpN"
$ws.Range("I2").Value = "a7q410000004I1W"
$ws.Range("J2").Value = "a7Z4100000000hb"
$ws.Range("K2").Value = "SY_ReceiptLoc"

# Row 3 data
$ws.Range("H3").Value = "a811K0000004fpN"
$ws.Range("I3").Value = "a7q410000004I1W"
$ws.Range("J3").Value = "a7Z4100000000hb"
$ws.Range("K3").Value = "SY_ReceiptLoc"

# Select entire column H downward to match resulting saved selection state
$ws.Range("A1:XFD1048576").Select()
$ws.Range("H1").Activate()

# Autofit columns to reflect resulting column widths
$ws.Columns.AutoFit() | Out-Null
